# BitWise Support EnumType AutoGenerate
# Adds two new columns (Properties / MonsterType) to the Example sheet,
# annotates the three header cells with explanatory comments, widens the
# new "Properties" column, and restores the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column data -------------------------------------------------
# Row 1 is the header row; rows 2-11 are the per-monster data rows that
# already exist in columns A:D (Id, Name, Strength, Difficulty).

$properties = @(
    "Properties",
    "不可阻挡",
    "牛逼|沉默",
    "晕眩",
    "晕眩",
    "晕眩",
    "晕眩",
    "晕眩|不可阻挡",
    "晕眩",
    "晕眩|八八八",
    "晕眩|da"
)

$monsterType = @(
    "MonsterType",
    "Monster",
    "Monster",
    "Monster",
    "Monster",
    "Monster",
    "Monster",
    "Humanoid",
    "Monster",
    "Humanoid",
    "Monster"
)

for ($i = 0; $i -lt $properties.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 5).Value = $properties[$i]
    $ws.Cells.Item($row, 6).Value = $monsterType[$i]
}

# --- Column width for the new "Properties" column (column E) ---------
$ws.Columns.Item(5).ColumnWidth = 13.93

# --- Header comments (enum / bitmask documentation) -------------------
$nl = [char]10

$dComment = "枚举" + $nl + "[Easy::1]" + $nl + "[Medium::2]" + $nl + "[Hard::3]"
$ws.Range("D1").AddComment($dComment) | Out-Null

$eComment = "位组合" + $nl + "[不可阻挡::1]" + $nl + "[晕眩::2]" + $nl + "[牛逼::3]" + $nl + "[沉默::4]" + $nl
$ws.Range("E1").AddComment($eComment) | Out-Null

$fComment = "枚举" + $nl + "[Humanoid::1]" + $nl + "[Monster::2]" + $nl + "[Npc::3]"
$ws.Range("F1").AddComment($fComment) | Out-Null

# --- Restore the active selection -------------------------------------
$ws.Range("I2").Select() | Out-Null
